# Opportunities & Engagement - 1st Merge - 12 Mar 2025
#
# Refresh the sample Opportunity/Engagement test data rows:
#  - Opportunities!A2/B2 get a new opportunity name/number, plus a new
#    opportunity-id value in C2.
#  - Engagements!A2/B2 get a new engagement name/number.

$wb = $excel.ActiveWorkbook

$opportunities = $wb.Worksheets.Item("Opportunities")
$opportunities.Range("A2").Value = "Project ACE - Lender Edu"
$opportunities.Range("B2").Value = "127133"
$opportunities.Range("C2").Value = 133775

$engagements = $wb.Worksheets.Item("Engagements")
$engagements.Range("A2").Value = "Project Apollo - Lender Education"
$engagements.Range("B2").Value = "124379"

# Restore the cursor positions Excel leaves behind when the sheets were
# last saved (Opportunities selection moved, Engagements stayed the
# active/selected tab with a new active cell).
[void]$opportunities.Range("G9").Select()

[void]$engagements.Activate()
[void]$engagements.Range("B6").Select()
